$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up the "high_speed" row: b/c values were too wide, tighten the range ---
$ws.Range("D42").Value = 125
$ws.Range("E42").Value = 125
$ws.Range("F42").Value = 150

# --- Add three new fuzzy sets for a "tone" variable (rows 43-45) ---
$ws.Range("A43").Value = "low_tone"
$ws.Range("B43").Value = "trapezoidal_mf"
$ws.Range("C43").Value = 70
$ws.Range("D43").Value = 90
$ws.Range("E43").Value = 90
$ws.Range("F43").Value = 100

$ws.Range("A44").Value = "mid_tone"
$ws.Range("B44").Value = "trapezoidal_mf"
$ws.Range("C44").Value = 80
$ws.Range("D44").Value = 100
$ws.Range("E44").Value = 100
$ws.Range("F44").Value = 120

$ws.Range("A45").Value = "high_tone"
$ws.Range("B45").Value = "trapezoidal_mf"
$ws.Range("C45").Value = 100
$ws.Range("D45").Value = 120
$ws.Range("E45").Value = 120
$ws.Range("F45").Value = 140

# --- Update the window scroll position / active selection to match the saved view ---
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H42").Select()
